$wb = $excel.ActiveWorkbook
$wsEnh = $wb.Worksheets.Item("Enhancements")
$wsQueries = $wb.Worksheets.Item("queries")

# --- Enhancements sheet: append two new rows (26 = blank yellow separator, 27 = new enhancement entry) ---

# Row 26: blank separator row styled like the other separator row (yellow fill + wrap text)
$sep = $wsEnh.Range("A26:H26")
$sep.Interior.Color = 65535
$sep.WrapText = $true

# Row 27: new enhancement entry
$wsEnh.Range("A27").Value = 42635
$wsEnh.Range("A27").NumberFormat = "mm-dd-yy"
$wsEnh.Range("A27").WrapText = $true

$wsEnh.Range("B27").Value = "mapping the settings data like the room types role types from the service call"
$wsEnh.Range("B27").WrapText = $true

$wsEnh.Range("C27").Value = "we already have service that gives all the settings information /user/settings(get) pls map the response of the room types fetch from the service call not binding the values in the controller if any"
$wsEnh.Range("C27").WrapText = $true

$wsEnh.Range("D27").Value = "rajashree"
$wsEnh.Range("D27").WrapText = $true

$wsEnh.Range("E27").Value = 42634
$wsEnh.Range("E27").NumberFormat = "mm-dd-yy"
$wsEnh.Range("E27").WrapText = $true

$wsEnh.Range("F27").Value = 42634
$wsEnh.Range("F27").NumberFormat = "mm-dd-yy"
$wsEnh.Range("F27").WrapText = $true

$wsEnh.Rows.Item(27).RowHeight = 75

# --- view/selection updates ---
# queries sheet selection moves to C2 (without leaving it the active tab)
$wsQueries.Range("C2").Select()
$wsEnh.Activate()

# Enhancements sheet: scroll the frozen pane and move the selection to C9
$wsEnh.Range("A17").Select()
$wsEnh.Range("C9").Select()
